$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Re-split the runs of the "Discretionary Access control (DAC) : " line.
#    The visible text is unchanged; Word's proofing pass broke the run that
#    used to read "Discretionary Access control (DAC)" into
#    "...(DAC" + ")" (both keeping the Heading2Char run style), and the
#    trailing " : " run into " :" + " ".  We reproduce the same run
#    boundaries by nudging a zero-effect character property on the
#    sub-ranges, which forces Word to split the backing runs without
#    altering any visible formatting.
# ---------------------------------------------------------------------------
$dacPara = $d.Paragraphs.Item(6)
$dacStart = $dacPara.Range.Start

$closeParen = $d.Range($dacStart + 33, $dacStart + 34)
$closeParen.Bold = 1
$closeParen.Bold = 0

$colonPart = $d.Range($dacStart + 34, $dacStart + 36)
$colonPart.Bold = 1
$colonPart.Bold = 0

# ---------------------------------------------------------------------------
# 2) Append the new content after "you create it you own it":
#      - a blank List Paragraph spacer
#      - Heading2 "Warehouse" + two List Paragraph bullets
#      - Heading2 "Metadata" + a closing Normal paragraph
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Last

# Blank spacer paragraph (List Paragraph style, no text)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last

# Heading2 "Warehouse"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Heading 2"
$p.Range.Text = "Warehouse"

# List Paragraph: "It is compute in snowflake not storage"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "List Paragraph"
$p.Range.Text = "It is compute in snowflake not storage"
$pStart = $p.Range.Start
$fullText = "It is compute in snowflake not storage"
$wIdx = $fullText.IndexOf("compute")
$wordRange = $d.Range($pStart + $wIdx, $pStart + $wIdx + "compute".Length)
$wordRange.Bold = 1
$wordRange.Bold = 0

# List Paragraph: "Elastic, scale up is called multi cluster warehouse"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "List Paragraph"
$p.Range.Text = "Elastic, scale up is called multi cluster warehouse"

# Heading2 "Metadata"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Heading 2"
$p.Range.Text = "Metadata"

# Normal paragraph about information_schema
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$p.Range.Text = "Meta data of each database is stored in information_schema of that database"
$pStart = $p.Range.Start
$fullText2 = "Meta data of each database is stored in information_schema of that database"
$sIdx = $fullText2.IndexOf("information_schema")
$schemaRange = $d.Range($pStart + $sIdx, $pStart + $sIdx + "information_schema".Length)
$schemaRange.Bold = 1
$schemaRange.Bold = 0
